$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data for the "Unknown" store-and-forward flag value
$ws.Range("A4").Value = "U"
$ws.Range("B4").Value = "Unknown"

# Resize the table to include the new row
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:B4"))

# Match the reported active cell/selection after the edit
$ws.Range("B10").Select()
